# Updated main GSC export data.
#
# The GSC export window rolled forward by one day: the oldest day
# (2025-10-10, the first data row under the header on the "Chart" sheet)
# drops out of the report, and every subsequent row shifts up by one,
# so the table shrinks from 84 data rows (A1:D85) to 83 data rows
# (A1:D84). Deleting the entire row 2 reproduces exactly that shift -
# all dates/values below move up one row and keep their original
# pairing, and the now-unused last row (85) disappears.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$ws.Rows.Item(2).Delete()
